$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Rename "DRA" suite to "DRAIAM" in cell A12
$ws.Range("A12").Value = "DRAIAM"

# Update the active selection on the sheet from F20 to B20
$ws.Range("B20").Select()
